# Update imputed values produced by the RandomForest algorithm run.
# (commit message: "Update Name of Algo")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.09190000000001
$ws.Range("D6").Value = -7.8976
$ws.Range("A14").Value = -21.6365
$ws.Range("C15").Value = -13.92289999999999
$ws.Range("A16").Value = -22.15170000000002
$ws.Range("D18").Value = -8.729800000000001
$ws.Range("D19").Value = -8.390199999999998
$ws.Range("A21").Value = -20.39659999999999
$ws.Range("C21").Value = -12.4974
$ws.Range("C22").Value = -11.8222
$ws.Range("A23").Value = -20.41079999999999
$ws.Range("C24").Value = -13.33339999999999
$ws.Range("A25").Value = -21.61629999999999
$ws.Range("A26").Value = -21.30959999999997
$ws.Range("C27").Value = -12.90939999999999
$ws.Range("C28").Value = -13.10299999999999
$ws.Range("A29").Value = -20.86519999999997
$ws.Range("D35").Value = -8.510899999999994
$ws.Range("C36").Value = -13.23030000000001
$ws.Range("C39").Value = -13.27670000000001
$ws.Range("A40").Value = -20.28540000000001
$ws.Range("D44").Value = -7.420499999999997
$ws.Range("C45").Value = -13.49089999999999
$ws.Range("D47").Value = -7.208800000000008
$ws.Range("C48").Value = -11.69789999999999
$ws.Range("C49").Value = -13.8508
$ws.Range("D50").Value = -8.386999999999999
$ws.Range("D51").Value = -8.124799999999997
$ws.Range("C52").Value = -11.0043
$ws.Range("D52").Value = -7.088199999999997
$ws.Range("A53").Value = -21.77219999999998
$ws.Range("C53").Value = -11.3361
$ws.Range("C54").Value = -13.1581
$ws.Range("D55").Value = -8.822000000000001
$ws.Range("A57").Value = -22.00899999999998
$ws.Range("C57").Value = -12.9384
$ws.Range("D57").Value = -8.570199999999994
$ws.Range("D58").Value = -8.189300000000005
$ws.Range("A59").Value = -22.54360000000002
$ws.Range("D64").Value = -7.472199999999991
$ws.Range("A65").Value = -21.82469999999998
$ws.Range("D66").Value = -7.547099999999992
$ws.Range("A69").Value = -21.60749999999999
$ws.Range("C70").Value = -11.4693
$ws.Range("C71").Value = -11.3242
$ws.Range("A79").Value = -20.654
$ws.Range("D80").Value = -7.685200000000002
$ws.Range("A83").Value = -22.2473
$ws.Range("D83").Value = -8.355499999999999
$ws.Range("C86").Value = -14.302
$ws.Range("C87").Value = -12.9213
$ws.Range("C89").Value = -14.03989999999999
$ws.Range("A91").Value = -20.91289999999999
$ws.Range("D92").Value = -6.650899999999996
$ws.Range("A93").Value = -21.36940000000001
$ws.Range("D94").Value = -6.856300000000001
$ws.Range("D96").Value = -8.234000000000002
$ws.Range("D97").Value = -8.963399999999982
$ws.Range("A100").Value = -22.25220000000001
$ws.Range("C101").Value = -12.4343
$ws.Range("D101").Value = -8.517700000000003
$ws.Range("A103").Value = -21.56839999999999
